{"js": "// Replace the worksheet date and the 25 two-digit multiplication\n// problems with the new values from the next day's worksheet.\n// Each old value is unique in the document, so a plain text search +\n// full replace of each hit is unambiguous and order independent.\nconst replacements = [\n  [\"2024-09-30 Monday\", \"2024-10-01 Tuesday\"],\n  [\"57\u00d793=\", \"66\u00d729=\"],\n  [\"82\u00d725=\", \"58\u00d796=\"],\n  [\"53\u00d792=\", \"46\u00d775=\"],\n  [\"75\u00d786=\", \"43\u00d784=\"],\n  [\"11\u00d722=\", \"17\u00d737=\"],\n  [\"46\u00d785=\", \"48\u00d753=\"],\n  [\"89\u00d745=\", \"65\u00d774=\"],\n  [\"43\u00d718=\", \"27\u00d719=\"],\n  [\"63\u00d763=\", \"43\u00d738=\"],\n  [\"43\u00d724=\", \"36\u00d742=\"],\n  [\"64\u00d753=\", \"35\u00d718=\"],\n  [\"18\u00d770=\", \"86\u00d731=\"],\n  [\"63\u00d759=\", \"54\u00d763=\"],\n  [\"63\u00d762=\", \"12\u00d793=\"],\n  [\"69\u00d776=\", \"66\u00d750=\"],\n  [\"53\u00d745=\", \"63\u00d767=\"],\n  [\"86\u00d753=\", \"62\u00d766=\"],\n  [\"24\u00d760=\", \"88\u00d741=\"],\n  [\"97\u00d761=\", \"83\u00d749=\"],\n  [\"93\u00d730=\", \"13\u00d734=\"],\n  [\"65\u00d798=\", \"60\u00d763=\"],\n  [\"45\u00d751=\", \"17\u00d777=\"],\n  [\"38\u00d784=\", \"13\u00d748=\"],\n  [\"37\u00d772=\", \"24\u00d748=\"],\n  [\"12\u00d784=\", \"17\u00d751=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the worksheet date and the 25 two-digit multiplication\n# problems with the new values from the next day's worksheet.\n# Each old value is unique in the document, so Find/Replace All for\n# each pair is unambiguous and order independent.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-09-30 Monday\", \"2024-10-01 Tuesday\"),\n    @(\"57\u00d793=\", \"66\u00d729=\"),\n    @(\"82\u00d725=\", \"58\u00d796=\"),\n    @(\"53\u00d792=\", \"46\u00d775=\"),\n    @(\"75\u00d786=\", \"43\u00d784=\"),\n    @(\"11\u00d722=\", \"17\u00d737=\"),\n    @(\"46\u00d785=\", \"48\u00d753=\"),\n    @(\"89\u00d745=\", \"65\u00d774=\"),\n    @(\"43\u00d718=\", \"27\u00d719=\"),\n    @(\"63\u00d763=\", \"43\u00d738=\"),\n    @(\"43\u00d724=\", \"36\u00d742=\"),\n    @(\"64\u00d753=\", \"35\u00d718=\"),\n    @(\"18\u00d770=\", \"86\u00d731=\"),\n    @(\"63\u00d759=\", \"54\u00d763=\"),\n    @(\"63\u00d762=\", \"12\u00d793=\"),\n    @(\"69\u00d776=\", \"66\u00d750=\"),\n    @(\"53\u00d745=\", \"63\u00d767=\"),\n    @(\"86\u00d753=\", \"62\u00d766=\"),\n    @(\"24\u00d760=\", \"88\u00d741=\"),\n    @(\"97\u00d761=\", \"83\u00d749=\"),\n    @(\"93\u00d730=\", \"13\u00d734=\"),\n    @(\"65\u00d798=\", \"60\u00d763=\"),\n    @(\"45\u00d751=\", \"17\u00d777=\"),\n    @(\"38\u00d784=\", \"13\u00d748=\"),\n    @(\"37\u00d772=\", \"24\u00d748=\"),\n    @(\"12\u00d784=\", \"17\u00d751=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
